$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# ---------------------------------------------------------------------------
# 1. Update the selection shown on the sheet (D30 -> G11)
# ---------------------------------------------------------------------------
$ws.Range("G11").Select()

# ---------------------------------------------------------------------------
# 2. Add the missing border edges so the C:G block for rows 2-5 gets a full
#    grid (matching the look of the other test-case sheets, e.g. sheet "Лист6").
#    Row 2 needs a top edge, column G needs a left edge.
# ---------------------------------------------------------------------------
$ws.Range("C2:G2").Borders.Item(8).Weight = 4
$ws.Range("G2:G5").Borders.Item(7).Weight = 4

# ---------------------------------------------------------------------------
# 3. Fill in the testcase content for row 2 (top half) and row 4 (bottom half)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "1.Check if the function takes all needed inputs"
$ws.Range("D2").Value = "1.Function successfully takes all input"
$ws.Range("E2").Value = "Passed"
$ws.Range("F2").Value = "1.Function takes the given inputs"
$ws.Range("G2").Value = "The test has`nbeen done manually"

$ws.Range("C4").Value = "2.Check if the function calculates everything right and returns a result"
$ws.Range("D4").Value = "2.The function calculates and returns a result"
$ws.Range("F4").Value = "2.The function calculates and returns a result"

# ---------------------------------------------------------------------------
# 4. Apply alignment / wrap formatting
# ---------------------------------------------------------------------------
$ws.Range("C2:C3").HorizontalAlignment = -4108
$ws.Range("C2:C3").WrapText = $true

$ws.Range("C4:C5").HorizontalAlignment = -4108
$ws.Range("C4:C5").VerticalAlignment = -4108
$ws.Range("C4:C5").WrapText = $true

$ws.Range("D2:D5").HorizontalAlignment = -4108
$ws.Range("D2:D5").VerticalAlignment = -4108
$ws.Range("D2:D3").WrapText = $true
$ws.Range("D4:D5").WrapText = $true

$ws.Range("E2:E5").HorizontalAlignment = -4108
$ws.Range("E2:E5").VerticalAlignment = -4108

$ws.Range("F2:F5").HorizontalAlignment = -4108
$ws.Range("F2:F5").VerticalAlignment = -4108
$ws.Range("F2:F3").WrapText = $true
$ws.Range("F4:F5").WrapText = $true

$ws.Range("G2:G5").HorizontalAlignment = -4108
$ws.Range("G2:G5").VerticalAlignment = -4108
$ws.Range("G2").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Merge the cells as in the other testcase sheets
# ---------------------------------------------------------------------------
$ws.Range("C2:C3").Merge()
$ws.Range("C4:C5").Merge()
$ws.Range("D2:D3").Merge()
$ws.Range("D4:D5").Merge()
$ws.Range("E2:E5").Merge()
$ws.Range("F2:F3").Merge()
$ws.Range("F4:F5").Merge()
$ws.Range("G2:G5").Merge()

# ---------------------------------------------------------------------------
# 6. Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 33.75
$ws.Rows.Item(4).RowHeight = 46.5
